# Refresh the crypto price/volume table (values scraped from coinranking.com).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.778.81"
$ws.Range("E2").Value = "  +1.30%  "

# Row 3
$ws.Range("D3").Value = "2.266.13"
$ws.Range("E3").Value = "  +0.75%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'303.53"
$ws.Range("E5").Value = "  +0.76%  "

# Row 6
$ws.Range("D6").Value = "'91.92"
$ws.Range("E6").Value = "  +0.94%  "

# Row 7
$ws.Range("E7").Value = "  +2.28%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").Value = "'0.483"
$ws.Range("E9").Value = "  +0.06%  "

# Row 10
$ws.Range("D10").Value = "'32.40"
$ws.Range("E10").Value = "  +2.07%  "

# Row 11
$ws.Range("D11").Value = "'53.41"
$ws.Range("E11").Value = "  +0.48%  "

# Row 12
$ws.Range("E12").Value = "  +0.54%  "

# Row 13
$ws.Range("E13").Value = "  -0.59%  "

# Row 14
$ws.Range("D14").Value = "'6.66"
$ws.Range("E14").Value = "  +1.43%  "

# Row 15
$ws.Range("D15").Value = "2.617.08"
$ws.Range("E15").Value = "  +0.68%  "

# Row 16
$ws.Range("D16").Value = "'14.22"
$ws.Range("E16").Value = "  +0.99%  "

# Row 17
$ws.Range("D17").Value = "2.287.02"
$ws.Range("E17").Value = "  +2.15%  "

# Row 18
$ws.Range("E18").Value = "  +2.64%  "

# Row 19
$ws.Range("D19").Value = "41.685.17"
$ws.Range("E19").Value = "  +1.23%  "

# Row 20
$ws.Range("D20").Value = "'12.46"
$ws.Range("E20").Value = "  +5.12%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0905"
$ws.Range("E21").Value = "  +0.49%  "

# Row 22
$ws.Range("E22").Value = "  +1.54%  "

# Row 23
$ws.Range("D23").Value = "'67.12"
$ws.Range("E23").Value = "  +0.66%  "

# Row 24
$ws.Range("D24").Value = "'239.65"
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("E25").Value = "  +1.46%  "

# Row 26
$ws.Range("E26").Value = "  +0.10%  "

# Row 27
$ws.Range("D27").Value = "'1.92"
$ws.Range("E27").Value = "  +3.51%  "

# Row 28
$ws.Range("D28").Value = "'23.90"
$ws.Range("E28").Value = "  +0.94%  "

# Row 29
$ws.Range("D29").Value = "'9.51"
$ws.Range("E29").Value = "  -0.39%  "

# Row 30
$ws.Range("E30").Value = "  -4.64%  "

# Row 31
$ws.Range("D31").Value = "'34.99"
$ws.Range("E31").Value = "  +4.86%  "

# Row 32
$ws.Range("D32").Value = "'160.90"
$ws.Range("E32").Value = "  +1.13%  "

# Row 33
$ws.Range("D33").Value = "'5.25"
$ws.Range("E33").Value = "  +1.65%  "

# Row 34
$ws.Range("E34").Value = "  -0.13%  "

# Row 35
$ws.Range("D35").Value = "'0.0743"
$ws.Range("E35").Value = "  +1.52%  "

# Row 36
$ws.Range("E36").Value = "  -0.15%  "

# Row 37
$ws.Range("D37").Value = "'16.90"
$ws.Range("E37").Value = "  +3.13%  "

# Row 38
$ws.Range("E38").Value = "  -0.04%  "

# Row 39
$ws.Range("E39").Value = "  +1.60%  "

# Row 40
$ws.Range("E40").Value = "  +0.06%  "

# Row 41
$ws.Range("E41").Value = "  +0.58%  "

# Row 42
$ws.Range("D42").Value = "'3.91"
$ws.Range("E42").Value = "  -0.26%  "

# Row 43
$ws.Range("D43").Value = "2.019.28"
$ws.Range("E43").Value = "  -2.81%  "

# Row 44
$ws.Range("D44").Value = "'19.22"
$ws.Range("E44").Value = "  -4.57%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0279"
$ws.Range("E45").Value = "  +0.60%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'10.27"
$ws.Range("E46").Value = "  +0.25%  "

# Row 47
$ws.Range("D47").Value = "'2.12"
$ws.Range("E47").Value = "  +4.31%  "

# Row 48
$ws.Range("D48").Value = "'2.88"
$ws.Range("E48").Value = "  -3.11%  "

# Row 49
$ws.Range("E49").Value = "  +0.90%  "

# Row 50
$ws.Range("E50").Value = "  +0.70%  "

# Row 51
$ws.Range("D51").Value = "'52.22"
$ws.Range("E51").Value = "  +3.11%  "
